# Updated symbol list on Tue Dec 20 15:51:35 UTC 2022 with GitHub Actions
#
# All "price"-like values in column D (and a couple of composite strings in
# column E) are stored as plain text in this workbook, even though they look
# like numbers. Writing a numeric-looking string straight into .Value makes
# Excel coerce it to a real number (losing formatting / trailing zeros), so
# for every such cell we momentarily force a Text number format, assign the
# literal string, then restore the cell's original ("Normal") style so no
# spurious formatting change is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $value) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# --- simple price-only updates -------------------------------------------
Set-TextValue "D2"  "250.96"
Set-TextValue "D3"  "22.99"
Set-TextValue "D4"  "5.498"
Set-TextValue "D5"  "0.05653"
Set-TextValue "D6"  "3.438"
Set-TextValue "D7"  "6.407"
Set-TextValue "D8"  "0.8209"
Set-TextValue "D9"  "0.9275"
Set-TextValue "D10" "0.1441"
Set-TextValue "D11" "0.07468"
Set-TextValue "D12" "0.03160"
Set-TextValue "D13" "0.03079"
Set-TextValue "D14" "0.09349"
Set-TextValue "D15" "3.558"
Set-TextValue "D16" "0.001608"
Set-TextValue "D17" "0.04735"

# Row 18 (One/ONE): price update + the "Worst in 24h" tag moves off this row
Set-TextValue "D18" "0.0005786"
$ws.Range("E18").Value = "17OneONE"

Set-TextValue "D19" "0.006363"
Set-TextValue "D20" "0.005027"
Set-TextValue "D21" "0.001031"
Set-TextValue "D22" "0.0001500"
Set-TextValue "D23" "3.731"
Set-TextValue "D24" "2.188"
Set-TextValue "D25" "0.3291"
Set-TextValue "D26" "0.1304"
Set-TextValue "D28" "0.0002998"
Set-TextValue "D40" "0.04016"

# --- rows 41-43: KickToken / BKEXToken / CEJI rotate positions -----------
# (row 41 was KickToken, row 42 was BKEXToken, row 43 was CEJI; they now
# appear in the order BKEXToken, CEJI, KickToken, each with refreshed data
# and the "Worst in 24h" tag landing on the new KickToken row)
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D41" "0.1075"
$ws.Range("E41").Value = "40BKEXTokenBKK"

$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D42" "0.002719"
$ws.Range("E42").Value = "41CEJICEJI"

$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue "D43" "0.002942"
$ws.Range("E43").Value = "42KickTokenKICKWorstin24h"

Set-TextValue "D44" "0.007608"
Set-TextValue "D45" "0.00005569"
Set-TextValue "D46" "0.00000000750"
Set-TextValue "D48" "0.6596"
Set-TextValue "D49" "0.2245"
Set-TextValue "D50" "0.00002099"
Set-TextValue "D51" "0.01009"
